$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.987.64'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.826.17'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9970'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6305'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9988'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07467'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2935'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.03'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07692'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('D12').Value = '1.829.28'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.983'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6667'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.97'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009590'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.047'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.25%  '
$ws.Range('D18').Value = '28.987.92'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '225.71'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9975'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.132'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9989'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '160.17'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1413'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.487'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.495'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.123'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.053'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05427'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.198'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.851'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7425'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.134'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.627'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.31%  '
$ws.Range('D37').Value = '1.238.56'
$ws.Range('E37').Value = '  -2.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.744'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01774'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.64%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.643'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8977'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9990'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.26'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('D44').Value = '1.976.94'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('E45').Value = '  +0.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.07'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5088'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4043'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.69%  '
$ws.Range('B49').Value = 'XinFinNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07263'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.65%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.904'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.658'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.24%  '
